# Excel COM-interop script that applies the scheduled market-data refresh
# described by the Ultima_Profits.xlsx diff ("chore: update Sheets via
# scheduled runner"). Each FFXIV Leve-profit sheet (ALC, ARM, BSM, CRP,
# CUL, GSM, LTW, WVR) has a handful of rows whose market-price columns
# (H:N = currentAveragePrice/NQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ) were
# refreshed with new data. A few cells that previously held 0/duplicate
# placeholder values are newly populated, and a few previously-populated
# cells are cleared, matching the source diff exactly.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I62").Value = 1300
$ws.Range("K62").Value = 1300
$ws.Range("M62").Value = -676

$ws.Range("I65").Value = 1300
$ws.Range("K65").Value = 6500
$ws.Range("M65").Value = -3380

$ws.Range("H137").Value = 7016.364
$ws.Range("I137").Value = 800.2308
$ws.Range("J137").Value = 11056.85
$ws.Range("K137").Value = 2400.6924
$ws.Range("L137").Value = 33170.55
$ws.Range("M137").Value = 149.3076000000001
$ws.Range("N137").Value = -38270.55

$ws.Range("H138").Value = 14498334
$ws.Range("I138").Value = 37041860
$ws.Range("J138").Value = 6067.857
$ws.Range("K138").Value = 111125580
$ws.Range("L138").Value = 18203.571
$ws.Range("M138").Value = -111120440
$ws.Range("N138").Value = -28483.571

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20756
$ws.Range("I32").Value = 20756
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 20756
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -20469
$ws.Range("N32").ClearContents()

$ws.Range("H45").Value = 1516680.6
$ws.Range("I45").Value = 2274210.5
$ws.Range("J45").Value = 1620.5
$ws.Range("K45").Value = 2274210.5
$ws.Range("L45").Value = 1620.5
$ws.Range("M45").Value = -2273833.5
$ws.Range("N45").Value = -2374.5

$ws.Range("H110").Value = 579.0714
$ws.Range("I110").Value = 557.6667
$ws.Range("J110").Value = 707.5
$ws.Range("K110").Value = 557.6667
$ws.Range("L110").Value = 707.5
$ws.Range("M110").Value = 1487.3333
$ws.Range("N110").Value = -4797.5

$ws.Range("H132").Value = 1121.1384
$ws.Range("I132").Value = 918.4286
$ws.Range("K132").Value = 2755.2858
$ws.Range("M132").Value = -225.2857999999997

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2109.5518
$ws.Range("I86").Value = 1955.7727
$ws.Range("K86").Value = 1955.7727
$ws.Range("M86").Value = -832.7727

$ws.Range("H89").Value = 2109.5518
$ws.Range("I89").Value = 1955.7727
$ws.Range("K89").Value = 9778.863499999999
$ws.Range("M89").Value = -4162.863499999999

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6042.8438
$ws.Range("I31").Value = 2013.2222
$ws.Range("J31").Value = 27802.8
$ws.Range("K31").Value = 2013.2222
$ws.Range("L31").Value = 27802.8
$ws.Range("M31").Value = -1718.2222
$ws.Range("N31").Value = -28392.8

$ws.Range("H34").Value = 6042.8438
$ws.Range("I34").Value = 2013.2222
$ws.Range("J34").Value = 27802.8
$ws.Range("K34").Value = 2013.2222
$ws.Range("L34").Value = 27802.8
$ws.Range("M34").Value = -1811.2222
$ws.Range("N34").Value = -28206.8

$ws.Range("H39").Value = 8015
$ws.Range("I39").Value = 8015
$ws.Range("K39").Value = 8015
$ws.Range("M39").Value = -7624

$ws.Range("H49").Value = 8015
$ws.Range("I49").Value = 8015
$ws.Range("K49").Value = 8015
$ws.Range("M49").Value = -7833

$ws.Range("H141").Value = 38926
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 38926
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 38926
$ws.Range("N141").Value = -49286
$ws.Range("M141").ClearContents()

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 1899.5
$ws.Range("J80").Value = 1900
$ws.Range("L80").Value = 5700
$ws.Range("N80").Value = -7572

$ws.Range("H83").Value = 1899.5
$ws.Range("J83").Value = 1900
$ws.Range("L83").Value = 17100
$ws.Range("N83").Value = -26460

$ws.Range("H107").Value = 438.38235
$ws.Range("I107").Value = 205.41667
$ws.Range("J107").Value = 565.4545000000001
$ws.Range("K107").Value = 616.25001
$ws.Range("L107").Value = 1696.3635
$ws.Range("M107").Value = 1303.74999
$ws.Range("N107").Value = -5536.3635

$ws.Range("H122").Value = 1242.6666
$ws.Range("I122").Value = 1258.1666
$ws.Range("J122").Value = 1222
$ws.Range("K122").Value = 11323.4994
$ws.Range("L122").Value = 10998
$ws.Range("M122").Value = -8873.499400000001
$ws.Range("N122").Value = -15898

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2680
$ws.Range("I80").Value = 2386.1428
$ws.Range("J80").Value = 2838.2307
$ws.Range("K80").Value = 2386.1428
$ws.Range("L80").Value = 2838.2307
$ws.Range("M80").Value = -1388.1428
$ws.Range("N80").Value = -4834.2307

$ws.Range("H83").Value = 2680
$ws.Range("I83").Value = 2386.1428
$ws.Range("J83").Value = 2838.2307
$ws.Range("K83").Value = 11930.714
$ws.Range("L83").Value = 14191.1535
$ws.Range("M83").Value = -6938.714
$ws.Range("N83").Value = -24175.1535

$ws.Range("H102").Value = 3161.35
$ws.Range("I102").Value = 4461.1665
$ws.Range("J102").Value = 1211.625
$ws.Range("K102").Value = 4461.1665
$ws.Range("L102").Value = 1211.625
$ws.Range("M102").Value = -2839.1665
$ws.Range("N102").Value = -4455.625

$ws.Range("H126").Value = 4292
$ws.Range("I126").Value = 3332.6667
$ws.Range("J126").Value = 4579.8
$ws.Range("K126").Value = 9998.000100000001
$ws.Range("L126").Value = 13739.4
$ws.Range("M126").Value = -7528.000100000001
$ws.Range("N126").Value = -18679.4

$ws.Range("H132").Value = 3409.8462
$ws.Range("I132").Value = 4228.7427
$ws.Range("J132").Value = 1723.8823
$ws.Range("K132").Value = 12686.2281
$ws.Range("L132").Value = 5171.6469
$ws.Range("M132").Value = -10156.2281
$ws.Range("N132").Value = -10231.6469

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2408.12
$ws.Range("I68").Value = 2193.75
$ws.Range("J68").Value = 2789.2222
$ws.Range("K68").Value = 2193.75
$ws.Range("L68").Value = 2789.2222
$ws.Range("M68").Value = -1444.75
$ws.Range("N68").Value = -4287.2222

$ws.Range("H71").Value = 2408.12
$ws.Range("I71").Value = 2193.75
$ws.Range("J71").Value = 2789.2222
$ws.Range("K71").Value = 10968.75
$ws.Range("L71").Value = 13946.111
$ws.Range("M71").Value = -7224.75
$ws.Range("N71").Value = -21434.111

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 26792.834
$ws.Range("J54").Value = 28213.4
$ws.Range("L54").Value = 28213.4
$ws.Range("N54").Value = -29253.4

$ws.Range("H81").Value = 540.9167
$ws.Range("I81").Value = 540.9167
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 1081.8334
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -20.83339999999998
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 540.9167
$ws.Range("I84").Value = 540.9167
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 5409.166999999999
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -105.1669999999995
$ws.Range("N84").ClearContents()

$ws.Range("H126").Value = 1899.6735
$ws.Range("I126").Value = 1890.174
$ws.Range("J126").Value = 1908.0769
$ws.Range("K126").Value = 5670.522
$ws.Range("L126").Value = 5724.2307
$ws.Range("M126").Value = -3200.522
$ws.Range("N126").Value = -10664.2307
